$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# 1. Collapse the stack of empty paragraphs (and the anchored picture they
#    held) that used to sit between the top bookmark paragraph and the
#    "Use Case Order Pizza" heading, leaving a single empty paragraph that
#    carries the _GoBack bookmark.
# -------------------------------------------------------------------------
$p2Start = $d.Paragraphs(2).Range.Start
$p20End  = $d.Paragraphs(20).Range.End
$gap = $d.Range($p2Start, $p20End)
$gap.Delete()

# Re-create the _GoBack bookmark (it lived inside the paragraph that was
# just removed) at the end of what is now the first, now-empty, paragraph.
$p1End = $d.Paragraphs(1).Range.End
$bmRange = $d.Range($p1End - 1, $p1End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# -------------------------------------------------------------------------
# 2. Drop the stale <w:lastRenderedPageBreak/> markers that Word had
#    cached in front of three of the use-case headings. Re-asserting the
#    same text through Find/Replace forces Word to rebuild those runs
#    without the cached page-break marker.
# -------------------------------------------------------------------------
$headings = @(
    "Use Case Order Pizza",
    "Use Case Choose Size and Topping",
    "5. Use Case Accept Bill"
)
foreach ($heading in $headings) {
    $d.Content.Find.Execute($heading, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $heading, 2) | Out-Null
}
